# Update the metric values on the active sheet to reflect the new
# vote-system (borda.count) results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - MAE
$ws.Range("B2").Value = 1.242
$ws.Range("C2").Value = 1.801
$ws.Range("D2").Value = 1.142
$ws.Range("E2").Value = 2.492
$ws.Range("F2").Value = 1.285

# Row 3 - MSE
$ws.Range("B3").Value = 2.928
$ws.Range("C3").Value = 6.29
$ws.Range("D3").Value = 2.175
$ws.Range("E3").Value = 14.225
$ws.Range("F3").Value = 3.344

# Row 4 - mean Y-Test
$ws.Range("B4").Value = 18.214
$ws.Range("C4").Value = 15.308
$ws.Range("D4").Value = 12.948
$ws.Range("E4").Value = 30.588
$ws.Range("F4").Value = 18.064

# Row 5 - mean Y-predicted
$ws.Range("B5").Value = 18.057
$ws.Range("C5").Value = 15.074
$ws.Range("D5").Value = 12.675
$ws.Range("E5").Value = 32.428
$ws.Range("F5").Value = 18.19

# Row 6 - R2
$ws.Range("B6").Value = 0.751
$ws.Range("C6").Value = 0.629
$ws.Range("D6").Value = 0.598
$ws.Range("E6").Value = 0.672
$ws.Range("F6").Value = 0.85
